$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update "Valor Mora" total (E11): 104000 -> 156000
$ws.Range("E11").Value = 156000

# 2) Update "Cant. Periodos" count (F13): 2 -> 3
$ws.Range("F13").Value = 3

# 3) Insert a new detail row above the existing second period row (old row 17),
#    so the table grows from 2 period rows to 3. The new row takes on the same
#    look/formatting as the row directly above it (row 16).
$ws.Rows("17:17").Insert()

$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Fill in the new row's data (same worker, same amounts, just a different period)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1006286414"
$ws.Range("D17").Value = "JHOAN ALBERTO ROMERO DIAZ"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 52000
$ws.Range("G17").Value = 1300000

# 4) Renumber the periods shown in the table: newest period on top.
$ws.Range("E16").Value = "2507"
$ws.Range("E18").Value = "2505"
